$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.051.14"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.960.86"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0851"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +78.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "3.424.06"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.14%  "
$ws.Range("E17").Value = "  +4.58%  "
$ws.Range("D18").Value = "2.961.11"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "51.018.69"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "0.0₃0953"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "266.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.165"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.88%  "
$ws.Range("E31").Value = "  -3.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.27%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0434"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.33%  "
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "2.011.14"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("E49").Value = "  -4.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0321"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.93%  "
